{"js": "// Replace the multiplication-problem text runs (e.g. \"163\u00d79=\") in the\n// document's single big table with their updated values, per the diff.\n// Old -> new values are all distinct, so a direct search/replace of each\n// exact old string is unambiguous and order-independent.\nconst replacements = [\n  [\"163\u00d79=\", \"444\u00d74=\"],\n  [\"144\u00d74=\", \"994\u00d75=\"],\n  [\"225\u00d78=\", \"713\u00d72=\"],\n  [\"485\u00d78=\", \"788\u00d75=\"],\n  [\"337\u00d79=\", \"779\u00d75=\"],\n  [\"699\u00d75=\", \"395\u00d79=\"],\n  [\"918\u00d73=\", \"896\u00d79=\"],\n  [\"930\u00d78=\", \"677\u00d79=\"],\n  [\"341\u00d76=\", \"843\u00d74=\"],\n  [\"856\u00d73=\", \"755\u00d76=\"],\n  [\"927\u00d78=\", \"456\u00d77=\"],\n  [\"486\u00d74=\", \"937\u00d77=\"],\n  [\"312\u00d75=\", \"587\u00d77=\"],\n  [\"558\u00d73=\", \"153\u00d77=\"],\n  [\"525\u00d78=\", \"428\u00d75=\"],\n  [\"346\u00d78=\", \"664\u00d72=\"],\n  [\"295\u00d74=\", \"956\u00d77=\"],\n  [\"226\u00d79=\", \"997\u00d72=\"],\n  [\"504\u00d79=\", \"910\u00d77=\"],\n  [\"334\u00d73=\", \"900\u00d77=\"],\n  [\"770\u00d75=\", \"854\u00d79=\"],\n  [\"826\u00d74=\", \"471\u00d78=\"],\n  [\"323\u00d77=\", \"350\u00d75=\"],\n  [\"593\u00d78=\", \"759\u00d74=\"],\n  [\"120\u00d75=\", \"489\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text runs (e.g. \"163\u00d79=\") in the\n# document's single big table with their updated values, per the diff.\n# Old -> new values are all distinct, so a direct Find/Replace of each\n# exact old string is unambiguous and order-independent.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n  @(\"163\u00d79=\", \"444\u00d74=\"),\n  @(\"144\u00d74=\", \"994\u00d75=\"),\n  @(\"225\u00d78=\", \"713\u00d72=\"),\n  @(\"485\u00d78=\", \"788\u00d75=\"),\n  @(\"337\u00d79=\", \"779\u00d75=\"),\n  @(\"699\u00d75=\", \"395\u00d79=\"),\n  @(\"918\u00d73=\", \"896\u00d79=\"),\n  @(\"930\u00d78=\", \"677\u00d79=\"),\n  @(\"341\u00d76=\", \"843\u00d74=\"),\n  @(\"856\u00d73=\", \"755\u00d76=\"),\n  @(\"927\u00d78=\", \"456\u00d77=\"),\n  @(\"486\u00d74=\", \"937\u00d77=\"),\n  @(\"312\u00d75=\", \"587\u00d77=\"),\n  @(\"558\u00d73=\", \"153\u00d77=\"),\n  @(\"525\u00d78=\", \"428\u00d75=\"),\n  @(\"346\u00d78=\", \"664\u00d72=\"),\n  @(\"295\u00d74=\", \"956\u00d77=\"),\n  @(\"226\u00d79=\", \"997\u00d72=\"),\n  @(\"504\u00d79=\", \"910\u00d77=\"),\n  @(\"334\u00d73=\", \"900\u00d77=\"),\n  @(\"770\u00d75=\", \"854\u00d79=\"),\n  @(\"826\u00d74=\", \"471\u00d78=\"),\n  @(\"323\u00d77=\", \"350\u00d75=\"),\n  @(\"593\u00d78=\", \"759\u00d74=\"),\n  @(\"120\u00d75=\", \"489\u00d75=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n}\n"}
